$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 14259.35
$ws.Range("I62").Value = 31327.285
$ws.Range("J62").Value = 5068.923
$ws.Range("K62").Value = 31327.285
$ws.Range("L62").Value = 5068.923
$ws.Range("M62").Value = -30703.285
$ws.Range("N62").Value = -6316.923
$ws.Range("H65").Value = 14259.35
$ws.Range("I65").Value = 31327.285
$ws.Range("J65").Value = 5068.923
$ws.Range("K65").Value = 156636.425
$ws.Range("L65").Value = 25344.615
$ws.Range("M65").Value = -153516.425
$ws.Range("N65").Value = -31584.615
$ws.Range("H98").Value = 1548.1471
$ws.Range("I98").Value = 928.88464
$ws.Range("J98").Value = 3560.75
$ws.Range("K98").Value = 928.88464
$ws.Range("L98").Value = 3560.75
$ws.Range("M98").Value = 569.11536
$ws.Range("N98").Value = -6556.75
$ws.Range("H121").Value = 1700
$ws.Range("I121").Value = 598.3333
$ws.Range("J121").Value = 3352.5
$ws.Range("K121").Value = 1794.9999
$ws.Range("L121").Value = 10057.5
$ws.Range("M121").Value = -47.99990000000003
$ws.Range("N121").Value = -13551.5
$ws.Range("H122").Value = 1548.1471
$ws.Range("I122").Value = 928.88464
$ws.Range("J122").Value = 3560.75
$ws.Range("K122").Value = 2786.65392
$ws.Range("L122").Value = 10682.25
$ws.Range("M122").Value = -336.6539199999997
$ws.Range("N122").Value = -15582.25
$ws.Range("H135").Value = 525.3261
$ws.Range("I135").Value = 361.07144
$ws.Range("J135").Value = 2250
$ws.Range("K135").Value = 3249.64296
$ws.Range("L135").Value = 20250
$ws.Range("M135").Value = -714.6429600000001
$ws.Range("N135").Value = -25320
$ws.Range("H137").Value = 4048.3704
$ws.Range("I137").Value = 4405.722
$ws.Range("J137").Value = 3333.6667
$ws.Range("K137").Value = 13217.166
$ws.Range("L137").Value = 10001.0001
$ws.Range("M137").Value = -10667.166
$ws.Range("N137").Value = -15101.0001
$ws.Range("H138").Value = 2851.5867
$ws.Range("I138").Value = 1576.4375
$ws.Range("J138").Value = 5118.5186
$ws.Range("K138").Value = 4729.3125
$ws.Range("L138").Value = 15355.5558
$ws.Range("M138").Value = 410.6875
$ws.Range("N138").Value = -25635.5558

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1084.5
$ws.Range("I45").Value = 980
$ws.Range("K45").Value = 980
$ws.Range("M45").Value = -603
$ws.Range("H74").Value = 270098.22
$ws.Range("I74").Value = 345865.97
$ws.Range("J74").Value = 86992.836
$ws.Range("K74").Value = 345865.97
$ws.Range("L74").Value = 86992.836
$ws.Range("M74").Value = -344991.97
$ws.Range("N74").Value = -88740.836
$ws.Range("H77").Value = 270098.22
$ws.Range("I77").Value = 345865.97
$ws.Range("J77").Value = 86992.836
$ws.Range("K77").Value = 1729329.85
$ws.Range("L77").Value = 434964.18
$ws.Range("M77").Value = -1724961.85
$ws.Range("N77").Value = -443700.18
$ws.Range("H88").Value = 4177.8887
$ws.Range("I88").Value = 6200.25
$ws.Range("J88").Value = 2560
$ws.Range("K88").Value = 6200.25
$ws.Range("L88").Value = 2560
$ws.Range("M88").Value = -5794.25
$ws.Range("N88").Value = -3372
$ws.Range("H91").Value = 4177.8887
$ws.Range("I91").Value = 6200.25
$ws.Range("J91").Value = 2560
$ws.Range("K91").Value = 6200.25
$ws.Range("L91").Value = 2560
$ws.Range("M91").Value = -4796.25
$ws.Range("N91").Value = -5368

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1513
$ws.Range("I94").Value = 1141.2727
$ws.Range("J94").Value = 2330.8
$ws.Range("K94").Value = 1141.2727
$ws.Range("L94").Value = 2330.8
$ws.Range("M94").Value = -690.2727
$ws.Range("N94").Value = -3232.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2454.3242
$ws.Range("I31").Value = 1880.0892
$ws.Range("J31").Value = 4240.8335
$ws.Range("K31").Value = 1880.0892
$ws.Range("L31").Value = 4240.8335
$ws.Range("M31").Value = -1585.0892
$ws.Range("N31").Value = -4830.8335
$ws.Range("H34").Value = 2454.3242
$ws.Range("I34").Value = 1880.0892
$ws.Range("J34").Value = 4240.8335
$ws.Range("K34").Value = 1880.0892
$ws.Range("L34").Value = 4240.8335
$ws.Range("M34").Value = -1678.0892
$ws.Range("N34").Value = -4644.8335
$ws.Range("H68").Value = 28326.666
$ws.Range("J68").Value = 28326.666
$ws.Range("L68").Value = 28326.666
$ws.Range("N68").Value = -29824.666
$ws.Range("H71").Value = 28326.666
$ws.Range("J71").Value = 28326.666
$ws.Range("L71").Value = 84979.99800000001
$ws.Range("N71").Value = -92467.99800000001
$ws.Range("H74").Value = 13896.25
$ws.Range("J74").Value = 13896.25
$ws.Range("L74").Value = 13896.25
$ws.Range("N74").Value = -15644.25
$ws.Range("H77").Value = 13896.25
$ws.Range("J77").Value = 13896.25
$ws.Range("L77").Value = 41688.75
$ws.Range("N77").Value = -50424.75
$ws.Range("H122").Value = 1522.95
$ws.Range("I122").Value = 944.6
$ws.Range("J122").Value = 2101.3
$ws.Range("K122").Value = 2833.8
$ws.Range("L122").Value = 6303.900000000001
$ws.Range("M122").Value = -383.8000000000002
$ws.Range("N122").Value = -11203.9
$ws.Range("H132").Value = 2084.1794
$ws.Range("I132").Value = 968.0741
$ws.Range("K132").Value = 2904.2223
$ws.Range("M132").Value = -374.2223000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1316.625
$ws.Range("I121").Value = 427.14285
$ws.Range("J121").Value = 1682.8823
$ws.Range("K121").Value = 1281.42855
$ws.Range("L121").Value = 5048.6469
$ws.Range("M121").Value = 28.57144999999991
$ws.Range("N121").Value = -7668.6469
$ws.Range("H137").Value = 3268.889
$ws.Range("I137").Value = 1782.5
$ws.Range("J137").Value = 3894.7368
$ws.Range("K137").Value = 5347.5
$ws.Range("L137").Value = 11684.2104
$ws.Range("M137").Value = -247.5
$ws.Range("N137").Value = -21884.2104

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4031.4614
$ws.Range("I102").Value = 2000.25
$ws.Range("J102").Value = 7281.4
$ws.Range("K102").Value = 2000.25
$ws.Range("L102").Value = 7281.4
$ws.Range("M102").Value = -378.25
$ws.Range("N102").Value = -10525.4
$ws.Range("H132").Value = 3353.9285
$ws.Range("I132").Value = 3188.7104
$ws.Range("K132").Value = 9566.1312
$ws.Range("M132").Value = -7036.1312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2801.6785
$ws.Range("I7").Value = 2741.0588
$ws.Range("K7").Value = 2741.0588
$ws.Range("M7").Value = -2629.0588
$ws.Range("H122").Value = 1791.5
$ws.Range("I122").Value = 1700
$ws.Range("J122").Value = 2249
$ws.Range("K122").Value = 5100
$ws.Range("L122").Value = 6747
$ws.Range("M122").Value = -2650
$ws.Range("N122").Value = -11647
$ws.Range("H126").Value = 2801.6785
$ws.Range("I126").Value = 2741.0588
$ws.Range("K126").Value = 8223.1764
$ws.Range("M126").Value = -5753.1764
$ws.Range("H132").Value = 5793.408
$ws.Range("I132").Value = 1578.3684
$ws.Range("J132").Value = 20354.455
$ws.Range("K132").Value = 4735.1052
$ws.Range("L132").Value = 61063.36500000001
$ws.Range("M132").Value = -2205.1052
$ws.Range("N132").Value = -66123.36500000001
$ws.Range("H136").Value = 3716.92
$ws.Range("I136").Value = 1995.8684
$ws.Range("J136").Value = 9166.916999999999
$ws.Range("K136").Value = 5987.6052
$ws.Range("L136").Value = 27500.751
$ws.Range("M136").Value = -3437.6052
$ws.Range("N136").Value = -32600.751
$ws.Range("H141").Value = 50357.5
$ws.Range("J141").Value = 50357.5
$ws.Range("L141").Value = 50357.5
$ws.Range("N141").Value = -60717.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 31282518
$ws.Range("I122").Value = 41709556
$ws.Range("J122").Value = 1406.875
$ws.Range("K122").Value = 125128668
$ws.Range("L122").Value = 4220.625
$ws.Range("M122").Value = -125126218
$ws.Range("N122").Value = -9120.625
$ws.Range("H126").Value = 1389.2273
$ws.Range("I126").Value = 855.3333
$ws.Range("J126").Value = 2029.9
$ws.Range("K126").Value = 2565.9999
$ws.Range("L126").Value = 6089.700000000001
$ws.Range("M126").Value = -95.9998999999998
$ws.Range("N126").Value = -11029.7
$ws.Range("H132").Value = 1766.4474
$ws.Range("I132").Value = 986.7143
$ws.Range("J132").Value = 3949.7
$ws.Range("K132").Value = 2960.1429
$ws.Range("L132").Value = 11849.1
$ws.Range("M132").Value = -430.1428999999998
$ws.Range("N132").Value = -16909.1
$ws.Range("H136").Value = 14043258
$ws.Range("I136").Value = 18538246
$ws.Range("J136").Value = 558294.75
$ws.Range("K136").Value = 55614738
$ws.Range("L136").Value = 1674884.25
$ws.Range("M136").Value = -55612188
$ws.Range("N136").Value = -1679984.25
